$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Window position (bookViews/workbookView xWindow/yWindow)
$win = $wb.Windows.Item(1)
$win.Left = 420
$win.Top = 1840

# Update comment texts in column C that changed content (order matches the
# order the author touched them in, so the shared-string table grows the same way)
$ws.Range("C8").Value = 'Saknar vinter/våraktivitet mellan 2011 och 2014. Har bara vinteraktivitet mellan 2000 och 2011 (vet ej vilka exakta datum). Har mailat Lars om 2011 -2014. Lars har mailat en del data tidigare år (c:a 2000 - 2005)  till Tomas Meijer och kanske till Anders. Alva kunde inte få ut något vettigt ur databasen. Det fattas årtal och datum på majoriteten av lybesöken. Det jag har nu är vinteraktivitet mellan 2000 och 2010 (BEBODDA_LYOR_HEF 00_10).'
$ws.Range("C14").Value = 'det finns rödrävsreproduktion i Peters fil fram till 2008 (röd text) dock är det bara två totalt. Använd koordinaterna för skjutna rävar istället. Finns fram till 2012 i Peters fil. (Helags_Red_Fox_Feeding).  Mailat Lars om resten. Använd Rasmus skript. '
$ws.Range("C17").Value = 'har mellan 2005 och 2008. Fick fler år av Karin också. Gå igenom.'
$ws.Range("C16").Value = 'Ingår fångstdatan i rastern jag fick av Rasmus? har fångstdata med lyprecision mellan 01 och 04 (Gnagfånst 2001-04 2004-09-16) och 08 och 14 (Sammanfattning08-14.xlsx). Fick en fil av Malin Larm om med en sammanfattning av gnagare 2008 -2017. 2009 och 2016 är dock inte med. Saknas alltså 2005 - 2007. Rasmus tror att datan för de åren kan ha dålig kvalitet. Fick en fil av Karin (Gnagarfångst2007) men den verkar bara ha för vindelfjällen för 2007.'
$ws.Range("C19").Value = ' 1,5 km eftersom rävar jagar närmare lyan när de har valpar (Frafjord 1993) och måste bära tillbaka mat till lyan (Zapata et al. 1998. Det är vettigare att ta en cirkelradie än riptrianglarna. Dessutom måste jag hålla observationerna oberoende av varandra. Gallant et al (2014) valde max radius på 1,5 km. '
$ws.Range("C20").Value = 'Den högsta medelvärdessannolikheten för lämmel i en lybuffer  är 0,53 under ett uppgångsår. Jag räknar allt från 0,265 och över som bra lämmelhabitat.  Medelbra är satt mellan 0,265 och medianen för medelvärdessannolikheten för lämmel i en lybuffer. Under medianen är dåliga  lämmelhabitat.'
$ws.Range("C9").Value = ' Gjort en sammanställning av kullar från fyra filer: "komplexa kullar Helags - genetiskt och observationer.csv"som har en sammanställning fram till 2017, Helags_Red_Fox_Feeding (från Peter, uppdaterad till 2008), BEBODDA_LYOR_HEF 00_10 (uppdaterad till 2010) och Fjellområder_kull_1977_2016_Bodil_Nina (som bara innehåller antalet kullar, inte specifika lyor. De stämmer inte överens med varandra. Sammanställningen heter "min sammanställning plus BEBODDA_LYOR_HEF 00_10.xlsx". Mailat Alva för att få ett utdrag ur databasen istället. Fick en ny fil av Bodil men den hjälpte inte. Typ som som Fjellområder_kull-filen.'

# Row 20 status changes from "påbörjat" to "klar" -> copy formatting from another
# "klar" cell (B11) so the font/color style matches, then set the new text
$ws.Range("B11").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B20").Value = "klar"

# Update selection to match the new active cell/range
$ws.Range("C12:C13").Select()

